$d = $word.ActiveDocument

# The document's single paragraph ends in the sentence
# "... la literatura universal." Someone placed the cursor right before
# the final period and typed in some gibberish text, so the sentence now
# reads "... la literatura universal nfsajibafjibadfuibafuibnsduiad."
#
# Locate the trailing period of the paragraph's sentence.
$para = $d.Paragraphs(1)
$sentenceRange = $para.Range
$sentenceRange.MoveEnd(1, -1) | Out-Null   # drop the paragraph mark
$end = $sentenceRange.End

# $end - 1 .. $end is exactly the final "." character.
$periodStart = $end - 1

# Insertion point sitting right before that trailing period.
$insertPoint = $d.Range($periodStart, $periodStart)
$insertPoint.InsertAfter(" nfsajibafjibadfuibafuibnsduiad")

# InsertAfter happily merges into the run it was inserted next to, but the
# source document models the inserted text and the trailing "." as their
# own runs. Nudge the freshly inserted span's direct formatting (flip a
# character property on and back off) so Word materialises it as separate
# run(s) with the same resulting formatting as their neighbours.
$insertedLen = " nfsajibafjibadfuibafuibnsduiad".Length
$insertedRange = $d.Range($periodStart, $periodStart + $insertedLen)
$insertedRange.Font.Bold = 1
$insertedRange.Font.Bold = 0
